$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.564.58'
$ws.Range("E2").Value = '  -4.42%  '
$ws.Range("D3").Value = '2.929.73'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''549.58'
$ws.Range("E5").Value = '  -4.30%  '
$ws.Range("D6").Value = '''130.23'
$ws.Range("E6").Value = '  +3.64%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.511'
$ws.Range("E8").Value = '  +1.54%  '
$ws.Range("D9").Value = '2.923.84'
$ws.Range("E9").Value = '  -2.45%  '
$ws.Range("D10").Value = '''0.127'
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("D11").Value = '''4.77'
$ws.Range("E11").Value = '  -5.64%  '
$ws.Range("D12").Value = '''0.446'
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("E13").Value = '  -0.03%  '
$ws.Range("D14").Value = '''32.88'
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("E15").Value = '  +0.25%  '
$ws.Range("D16").Value = '3.416.58'
$ws.Range("E16").Value = '  -2.36%  '
$ws.Range("D17").Value = '''6.86'
$ws.Range("E17").Value = '  +6.27%  '
$ws.Range("D18").Value = '2.927.89'
$ws.Range("E18").Value = '  -2.42%  '
$ws.Range("D19").Value = '57.617.44'
$ws.Range("E19").Value = '  -4.22%  '
$ws.Range("D20").Value = '''416.75'
$ws.Range("E20").Value = '  -2.90%  '
$ws.Range("D21").Value = '''13.18'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '''0.686'
$ws.Range("E22").Value = '  +2.58%  '
$ws.Range("D23").Value = '''6.98'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("D24").Value = '''13.03'
$ws.Range("E24").Value = '  +1.04%  '
$ws.Range("D25").Value = '''79.82'
$ws.Range("E25").Value = '  +0.51%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("E28").Value = '  -2.90%  '
$ws.Range("D29").Value = '''7.48'
$ws.Range("E29").Value = '  +2.99%  '
$ws.Range("E30").Value = '  +1.34%  '
$ws.Range("D31").Value = '''25.13'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").Value = '''6.00'
$ws.Range("E32").Value = '  -2.63%  '
$ws.Range("D33").Value = '''0.0968'
$ws.Range("E33").Value = '  +2.63%  '
$ws.Range("D34").Value = '''5.66'
$ws.Range("E34").Value = '  +0.91%  '
$ws.Range("D35").Value = '''0.942'
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = '''48.07'
$ws.Range("E37").Value = '  -4.48%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0683'
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").Value = '''8.70'
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("D40").Value = '''2.55'
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '''378.90'
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.107'
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("D43").Value = '''0.0346'
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").Value = '2.686.37'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D47").Value = '''121.92'
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").Value = '''1.98'
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("E50").Value = '  -2.23%  '
$ws.Range("E51").Value = '  -0.26%  '
